$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 10).Value = 4879
$ws.Cells.Item(3, 10).Value = 5158
$ws.Cells.Item(4, 2).Value = 1683
$ws.Cells.Item(4, 8).Value = 1700
$ws.Cells.Item(4, 10).Value = 1150
$ws.Cells.Item(5, 10).Value = 411
$ws.Cells.Item(6, 10).Value = 6397
$ws.Cells.Item(7, 2).Value = 23315
$ws.Cells.Item(7, 8).Value = 26011
$ws.Cells.Item(7, 10).Value = 17995

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 10).Value = 142
$ws.Cells.Item(7, 10).Value = 521
$ws.Cells.Item(8, 10).Value = 1146
$ws.Cells.Item(10, 10).Value = 121
$ws.Cells.Item(11, 10).Value = 276
$ws.Cells.Item(14, 10).Value = 79
$ws.Cells.Item(15, 10).Value = 193
$ws.Cells.Item(20, 10).Value = 376
$ws.Cells.Item(22, 10).Value = 49
$ws.Cells.Item(23, 10).Value = 174
$ws.Cells.Item(24, 10).Value = 51
$ws.Cells.Item(27, 10).Value = 98
$ws.Cells.Item(31, 10).Value = 164
$ws.Cells.Item(32, 10).Value = 29
$ws.Cells.Item(33, 10).Value = 816
$ws.Cells.Item(36, 10).Value = 251
$ws.Cells.Item(37, 10).Value = 559
$ws.Cells.Item(41, 10).Value = 116
$ws.Cells.Item(42, 10).Value = 733
$ws.Cells.Item(44, 10).Value = 135
$ws.Cells.Item(47, 10).Value = 139
$ws.Cells.Item(48, 10).Value = 201
$ws.Cells.Item(51, 10).Value = 226
$ws.Cells.Item(52, 10).Value = 452
$ws.Cells.Item(53, 10).Value = 233
$ws.Cells.Item(55, 10).Value = 228
$ws.Cells.Item(63, 2).Value = 389
$ws.Cells.Item(63, 8).Value = 255
$ws.Cells.Item(63, 10).Value = 72
$ws.Cells.Item(65, 10).Value = 473
$ws.Cells.Item(67, 10).Value = 697
$ws.Cells.Item(68, 10).Value = 34
$ws.Cells.Item(74, 10).Value = 20
$ws.Cells.Item(76, 10).Value = 260
$ws.Cells.Item(77, 10).Value = 141
$ws.Cells.Item(78, 10).Value = 228
$ws.Cells.Item(79, 10).Value = 516
$ws.Cells.Item(84, 10).Value = 150
$ws.Cells.Item(85, 10).Value = 787
$ws.Cells.Item(88, 10).Value = 197
$ws.Cells.Item(89, 10).Value = 228
$ws.Cells.Item(91, 10).Value = 200
$ws.Cells.Item(95, 10).Value = 270
$ws.Cells.Item(96, 10).Value = 216
$ws.Cells.Item(98, 10).Value = 114
$ws.Cells.Item(99, 10).Value = 273
$ws.Cells.Item(101, 2).Value = 23315
$ws.Cells.Item(101, 8).Value = 26011
$ws.Cells.Item(101, 10).Value = 17995

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 10).Value = 166
$ws.Cells.Item(3, 10).Value = 154
$ws.Cells.Item(4, 10).Value = 20
$ws.Cells.Item(7, 10).Value = 521

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(3, 10).Value = 59
$ws.Cells.Item(4, 10).Value = 27
$ws.Cells.Item(7, 10).Value = 228

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(2, 10).Value = 92
$ws.Cells.Item(4, 10).Value = 19
$ws.Cells.Item(6, 10).Value = 105
$ws.Cells.Item(7, 10).Value = 276

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(6, 10).Value = 191
$ws.Cells.Item(7, 10).Value = 452

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Cells.Item(3, 10).Value = 19
$ws.Cells.Item(6, 10).Value = 21
$ws.Cells.Item(7, 10).Value = 79

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(2, 10).Value = 67
$ws.Cells.Item(6, 10).Value = 74
$ws.Cells.Item(7, 10).Value = 216

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 10).Value = 323
$ws.Cells.Item(6, 10).Value = 384
$ws.Cells.Item(7, 10).Value = 1146

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(2, 10).Value = 47
$ws.Cells.Item(6, 10).Value = 144
$ws.Cells.Item(7, 10).Value = 233

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(3, 10).Value = 286
$ws.Cells.Item(6, 10).Value = 224
$ws.Cells.Item(7, 10).Value = 787

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 10).Value = 94
$ws.Cells.Item(7, 10).Value = 270

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(2, 10).Value = 78
$ws.Cells.Item(3, 10).Value = 103
$ws.Cells.Item(7, 10).Value = 273

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 10).Value = 173
$ws.Cells.Item(7, 10).Value = 697

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(2, 10).Value = 62
$ws.Cells.Item(7, 10).Value = 164

$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(3, 10).Value = 46
$ws.Cells.Item(7, 10).Value = 150

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(5, 10).Value = 24
$ws.Cells.Item(7, 10).Value = 559

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(3, 10).Value = 261
$ws.Cells.Item(4, 10).Value = 36
$ws.Cells.Item(7, 10).Value = 816

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(2, 10).Value = 134
$ws.Cells.Item(7, 10).Value = 473

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(3, 10).Value = 352
$ws.Cells.Item(4, 10).Value = 56

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(2, 10).Value = 44
$ws.Cells.Item(6, 10).Value = 46
$ws.Cells.Item(7, 10).Value = 135

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(2, 10).Value = 32
$ws.Cells.Item(7, 10).Value = 201

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(3, 10).Value = 54
$ws.Cells.Item(6, 10).Value = 139
$ws.Cells.Item(7, 10).Value = 260

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Cells.Item(2, 10).Value = 25
$ws.Cells.Item(7, 10).Value = 116

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(3, 10).Value = 148
$ws.Cells.Item(6, 10).Value = 375
$ws.Cells.Item(7, 10).Value = 733

$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(6, 10).Value = 68
$ws.Cells.Item(7, 10).Value = 121

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(3, 10).Value = 78
$ws.Cells.Item(6, 10).Value = 61
$ws.Cells.Item(7, 10).Value = 228

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(2, 10).Value = 54
$ws.Cells.Item(6, 10).Value = 110
$ws.Cells.Item(7, 10).Value = 228

$ws = $wb.Worksheets.Item('Dunning')
$ws.Cells.Item(3, 10).Value = 18
$ws.Cells.Item(7, 10).Value = 51

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(2, 10).Value = 46
$ws.Cells.Item(5, 10).Value = 6
$ws.Cells.Item(7, 10).Value = 174

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(6, 10).Value = 40
$ws.Cells.Item(7, 10).Value = 200

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 10).Value = 147
$ws.Cells.Item(6, 10).Value = 142
$ws.Cells.Item(7, 10).Value = 516

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 10).Value = 103
$ws.Cells.Item(7, 10).Value = 376

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(3, 10).Value = 78
$ws.Cells.Item(7, 10).Value = 251

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(6, 10).Value = 61
$ws.Cells.Item(7, 10).Value = 139

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(4, 10).Value = 7
$ws.Cells.Item(7, 10).Value = 193

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Cells.Item(4, 10).Value = 7
$ws.Cells.Item(6, 10).Value = 66
$ws.Cells.Item(7, 10).Value = 114

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(2, 10).Value = 41
$ws.Cells.Item(7, 10).Value = 142

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(2, 10).Value = 42
$ws.Cells.Item(7, 10).Value = 197

$ws = $wb.Worksheets.Item('Galewood')
$ws.Cells.Item(6, 10).Value = 13
$ws.Cells.Item(7, 10).Value = 29

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(2, 10).Value = 28
$ws.Cells.Item(7, 10).Value = 98

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(4, 10).Value = 24
$ws.Cells.Item(6, 10).Value = 80
$ws.Cells.Item(7, 10).Value = 226

$ws = $wb.Worksheets.Item('North Park')
$ws.Cells.Item(3, 10).Value = 9
$ws.Cells.Item(7, 10).Value = 34

$ws = $wb.Worksheets.Item('Clearing')
$ws.Cells.Item(3, 10).Value = 15
$ws.Cells.Item(7, 10).Value = 49

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Cells.Item(2, 10).Value = 49
$ws.Cells.Item(3, 10).Value = 51
$ws.Cells.Item(4, 10).Value = 14
$ws.Cells.Item(7, 10).Value = 141

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Cells.Item(6, 10).Value = 10
$ws.Cells.Item(7, 10).Value = 20

Write-Output "Applied 170 cell updates across 46 sheets"